{"js": "// Auto-generated edit script for Office.js (Word JavaScript API).\n// Applies the renumbering of [[PERSON_N]] placeholders across the\n// \"2. SEZNAM SUBJEKTU\" list and all later references in the document,\n// plus inserts 4 new list entries and removes 4 obsolete ones.\n\nconst DATA = {\n  \"replacements\": [\n    [\n      \"[[PERSON_5]] \u2013 \u201eo [[PERSON_6]]\u201c\",\n      \"[[PERSON_5]] \u2013 \u201eo [[PERSON_5]]\u201c\"\n    ],\n    [\n      \"[[PERSON_7]] \u2013 \u201ek [[PERSON_8]]\u201c\",\n      \"[[PERSON_6]] \u2013 \u201ek [[PERSON_7]]\u201c\"\n    ],\n    [\n      \"[[PERSON_9]] \u2013 \u201epro [[PERSON_10]]\u201c\",\n      \"[[PERSON_8]] \u2013 \u201epro [[PERSON_9]]\u201c\"\n    ],\n    [\n      \"[[PERSON_11]] \u2013 \u201es [[PERSON_12]]\u201c\",\n      \"[[PERSON_10]] \u2013 \u201es [[PERSON_10]]\u201c\"\n    ],\n    [\n      \"[[PERSON_13]] \u2013 \u201eu [[PERSON_13]]\u201c\",\n      \"[[PERSON_11]] \u2013 \u201eu [[PERSON_11]]\u201c\"\n    ],\n    [\n      \"[[PERSON_14]] \u2013 \u201eod [[PERSON_14]]\u201c\",\n      \"[[PERSON_12]] \u2013 \u201eod [[PERSON_12]]\u201c\"\n    ],\n    [\n      \"[[PERSON_15]] \u2013 \u201epro [[PERSON_15]]\u201c\",\n      \"[[PERSON_13]] \u2013 \u201epro [[PERSON_13]]\u201c\"\n    ],\n    [\n      \"[[PERSON_16]] \u2013 \u201ek [[PERSON_17]]\u201c\",\n      \"[[PERSON_14]] \u2013 \u201ek [[PERSON_15]]\u201c\"\n    ],\n    [\n      \"[[PERSON_18]] \u2013 \u201es [[PERSON_18]]\u201c\",\n      \"[[PERSON_16]] \u2013 \u201es [[PERSON_16]]\u201c\"\n    ],\n    [\n      \"[[PERSON_19]] \u2013 \u201eo [[PERSON_20]]\u201c\",\n      \"[[PERSON_17]] \u2013 \u201eo [[PERSON_18]]\u201c\"\n    ],\n    [\n      \"[[PERSON_21]] \u2013 \u201ek [[PERSON_22]]\u201c\",\n      \"[[PERSON_19]] \u2013 \u201ek [[PERSON_19]]\u201c\"\n    ],\n    [\n      \"[[PERSON_23]] \u2013 \u201es [[PERSON_24]]\u201c\",\n      \"[[PERSON_20]] \u2013 \u201es [[PERSON_20]]\u201c\"\n    ],\n    [\n      \"[[PERSON_26]] \u2013 \u201es [[PERSON_26]]\u201c\",\n      \"[[PERSON_26]] \u2013 \u201epro [[PERSON_27]]\u201c\"\n    ],\n    [\n      \"[[PERSON_27]] \u2013 \u201eo [[PERSON_27]]\u201c\",\n      \"[[PERSON_28]] \u2013 \u201es [[PERSON_28]]\u201c\"\n    ],\n    [\n      \"[[PERSON_28]] \u2013 \u201ek [[PERSON_28]]\u201c\",\n      \"[[PERSON_29]] \u2013 \u201ek [[PERSON_29]]\u201c\"\n    ],\n    [\n      \"[[PERSON_29]] \u2013 \u201eu [[PERSON_29]]\u201c\",\n      \"[[PERSON_30]] \u2013 \u201es [[PERSON_31]]\u201c\"\n    ],\n    [\n      \"[[PERSON_30]] \u2013 \u201epro [[PERSON_31]]\u201c\",\n      \"[[PERSON_32]] \u2013 \u201eo [[PERSON_33]]\u201c\"\n    ],\n    [\n      \"[[PERSON_32]] \u2013 \u201es [[PERSON_32]]\u201c\",\n      \"[[PERSON_34]] \u2013 \u201epro [[PERSON_34]]\u201c\"\n    ],\n    [\n      \"[[PERSON_33]] \u2013 \u201ek [[PERSON_33]]\u201c\",\n      \"[[PERSON_35]] \u2013 \u201es [[PERSON_36]]\u201c\"\n    ],\n    [\n      \"[[PERSON_34]] \u2013 \u201es [[PERSON_35]]\u201c\",\n      \"[[PERSON_37]] \u2013 \u201ek [[PERSON_38]]\u201c\"\n    ],\n    [\n      \"[[PERSON_36]] \u2013 \u201eo [[PERSON_37]]\u201c\",\n      \"[[PERSON_39]] \u2013 \u201es [[PERSON_39]]\u201c\"\n    ],\n    [\n      \"[[PERSON_38]] \u2013 \u201epro [[PERSON_38]]\u201c\",\n      \"[[PERSON_40]] \u2013 \u201eo [[PERSON_40]]\u201c\"\n    ],\n    [\n      \"[[PERSON_39]] \u2013 \u201es [[PERSON_40]]\u201c\",\n      \"[[PERSON_41]] \u2013 \u201eu [[PERSON_41]]\u201c\"\n    ],\n    [\n      \"[[PERSON_41]] \u2013 \u201ek [[PERSON_42]]\u201c\",\n      \"[[PERSON_42]] \u2013 \u201ek [[PERSON_42]]\u201c\"\n    ],\n    [\n      \"[[PERSON_43]] \u2013 \u201es [[PERSON_43]]\u201c\",\n      \"[[PERSON_43]] \u2013 \u201ese [[PERSON_44]]\u201c\"\n    ],\n    [\n      \"[[PERSON_44]] \u2013 \u201eo [[PERSON_45]]\u201c\",\n      \"[[PERSON_45]] \u2013 \u201eu [[PERSON_45]]\u201c\"\n    ],\n    [\n      \"[[PERSON_46]] \u2013 \u201eu [[PERSON_46]]\u201c\",\n      \"[[PERSON_46]] \u2013 \u201eo [[PERSON_47]]\u201c\"\n    ],\n    [\n      \"[[PERSON_47]] \u2013 \u201ek [[PERSON_47]]\u201c\",\n      \"[[PERSON_48]] \u2013 \u201es [[PERSON_48]]\u201c\"\n    ],\n    [\n      \"[[PERSON_48]] \u2013 \u201ese [[PERSON_49]]\u201c\",\n      \"[[PERSON_49]] \u2013 \u201ek [[PERSON_50]]\u201c\"\n    ],\n    [\n      \"[[PERSON_50]] \u2013 \u201eu [[PERSON_50]]\u201c\",\n      \"[[PERSON_51]] \u2013 \u201eod [[PERSON_51]]\u201c\"\n    ],\n    [\n      \"[[PERSON_51]] \u2013 \u201eo [[PERSON_52]]\u201c\",\n      \"[[PERSON_52]] \u2013 \u201es [[PERSON_52]]\u201c\"\n    ],\n    [\n      \"[[PERSON_53]] \u2013 \u201es [[PERSON_53]]\u201c\",\n      \"[[PERSON_53]] \u2013 \u201eu [[PERSON_53]]\u201c\"\n    ],\n    [\n      \"[[PERSON_54]] \u2013 \u201ek [[PERSON_55]]\u201c\",\n      \"[[PERSON_54]] \u2013 \u201eo [[PERSON_55]]\u201c\"\n    ],\n    [\n      \"[[PERSON_61]] \u2013 \u201ek [[PERSON_61]]\u201c\",\n      \"[[PERSON_56]] \u2013 \u201ek [[PERSON_56]]\u201c\"\n    ],\n    [\n      \"V t\u011bchto \u0159\u00edzen\u00edch bylo jedn\u00e1no nap\u0159. s [[PERSON_3]], [[PERSON_9]], [[PERSON_36]] \u010di [[PERSON_62]].\",\n      \"V t\u011bchto \u0159\u00edzen\u00edch bylo jedn\u00e1no nap\u0159. s [[PERSON_3]], [[PERSON_8]], [[PERSON_32]] \u010di [[PERSON_57]].\"\n    ],\n    [\n      \"sv\u011bdek [[PERSON_43]] (ve v\u00fdpov\u011bdi ozna\u010den jako \u201esv\u011bdek \u010cern\u00e9ho\u201c),\",\n      \"sv\u011bdek [[PERSON_39]] (ve v\u00fdpov\u011bdi ozna\u010den jako \u201esv\u011bdek \u010cern\u00e9ho\u201c),\"\n    ],\n    [\n      \"po\u0161kozen\u00e1 [[PERSON_21]] (\u201evyj\u00e1d\u0159en\u00ed [[PERSON_21]]\u201c),\",\n      \"po\u0161kozen\u00e1 [[PERSON_19]] (\u201evyj\u00e1d\u0159en\u00ed [[PERSON_19]]\u201c),\"\n    ],\n    [\n      \"ob\u017ealovan\u00fd [[PERSON_7]] (\u201eobhajoba [[PERSON_7]]\u201c),\",\n      \"ob\u017ealovan\u00fd [[PERSON_6]] (\u201eobhajoba [[PERSON_6]]\u201c),\"\n    ],\n    [\n      \"pr\u00e1vn\u00ed z\u00e1stupkyn\u011b JUDr. [[PERSON_47]], advok\u00e1tka,\",\n      \"pr\u00e1vn\u00ed z\u00e1stupkyn\u011b JUDr. [[PERSON_42]], advok\u00e1tka,\"\n    ],\n    [\n      \"tlumo\u010dn\u00edk [[PERSON_46]], zapsan\u00fd v seznamu tlumo\u010dn\u00edk\u016f.\",\n      \"tlumo\u010dn\u00edk [[PERSON_41]], zapsan\u00fd v seznamu tlumo\u010dn\u00edk\u016f.\"\n    ],\n    [\n      \"Alergologick\u00e9 vy\u0161et\u0159en\u00ed \u010d. ALG/2025/22751 proveden\u00e9 MUDr. [[PERSON_18]],\",\n      \"Alergologick\u00e9 vy\u0161et\u0159en\u00ed \u010d. ALG/2025/22751 proveden\u00e9 MUDr. [[PERSON_16]],\"\n    ],\n    [\n      \"Neurologick\u00e9 testy \u010d. NEU/2025/44119 proveden\u00e9 MUDr. [[PERSON_51]],\",\n      \"Neurologick\u00e9 testy \u010d. NEU/2025/44119 proveden\u00e9 MUDr. [[PERSON_46]],\"\n    ],\n    [\n      \"O\u010dn\u00ed vy\u0161et\u0159en\u00ed \u010d. OFT/2023/11281 proveden\u00e9 MUDr. [[PERSON_44]].\",\n      \"O\u010dn\u00ed vy\u0161et\u0159en\u00ed \u010d. OFT/2023/11281 proveden\u00e9 MUDr. [[PERSON_40]].\"\n    ],\n    [\n      \"Zvl\u00e1\u0161tn\u00ed pozornost byla v\u011bnov\u00e1na v\u00fdsledk\u016fm [[PERSON_25]], [[PERSON_29]] a [[PERSON_59]].\",\n      \"Zvl\u00e1\u0161tn\u00ed pozornost byla v\u011bnov\u00e1na v\u00fdsledk\u016fm [[PERSON_21]], [[PERSON_25]] a [[PERSON_54]].\"\n    ],\n    [\n      \"mobil [[PERSON_63]] S22, [[IMEI_1]],\",\n      \"mobil [[PERSON_58]] S22, [[IMEI_1]],\"\n    ],\n    [\n      \"[[PERSON_64]] poskytly technick\u00e9 p\u0159\u00edstupy pro \u0159e\u0161en\u00ed kauz:\",\n      \"[[PERSON_59]] poskytly technick\u00e9 p\u0159\u00edstupy pro \u0159e\u0161en\u00ed kauz:\"\n    ],\n    [\n      \"pr\u00e1vn\u00ed cloud \u00fa\u010det ID: LEX-ACC-88221 (spravovala [[PERSON_61]]),\",\n      \"pr\u00e1vn\u00ed cloud \u00fa\u010det ID: LEX-ACC-88221 (spravovala [[PERSON_56]]),\"\n    ],\n    [\n      \"[[PERSON_57]] (\u201ev\u00fdslech [[PERSON_57]]\u201c),\",\n      \"[[PERSON_52]] (\u201ev\u00fdslech [[PERSON_52]]\u201c),\"\n    ],\n    [\n      \"[[PERSON_53]] (\u201ev\u00fdpov\u011b\u010f [[PERSON_53]]\u201c),\",\n      \"[[PERSON_48]] (\u201ev\u00fdpov\u011b\u010f [[PERSON_48]]\u201c),\"\n    ],\n    [\n      \"[[PERSON_50]] (\u201ez\u00e1znam o v\u00fdslechu [[PERSON_50]]\u201c),\",\n      \"[[PERSON_45]] (\u201ez\u00e1znam o v\u00fdslechu [[PERSON_45]]\u201c),\"\n    ],\n    [\n      \"[[PERSON_28]] (\u201ev\u00fdslech [[PERSON_65]]\u201c).\",\n      \"[[PERSON_24]] (\u201ev\u00fdslech [[PERSON_60]]\u201c).\"\n    ],\n    [\n      \"PhDr. [[PERSON_44]] \u2013 psychologick\u00fd posudek,\",\n      \"PhDr. [[PERSON_40]] \u2013 psychologick\u00fd posudek,\"\n    ],\n    [\n      \"MUDr. [[PERSON_36]] \u2013 posudek z traumatologie,\",\n      \"MUDr. [[PERSON_32]] \u2013 posudek z traumatologie,\"\n    ],\n    [\n      \"Ing. [[PERSON_14]] \u2013 expertiza IT infrastruktury.\",\n      \"Ing. [[PERSON_12]] \u2013 expertiza IT infrastruktury.\"\n    ],\n    [\n      \"Tyto \u00fa\u010dty byly dolo\u017eeny nap\u0159. od [[PERSON_30]], [[PERSON_54]] nebo [[PERSON_66]].\",\n      \"Tyto \u00fa\u010dty byly dolo\u017eeny nap\u0159. od [[PERSON_26]], [[PERSON_49]] nebo [[PERSON_61]].\"\n    ],\n    [\n      \"[[PERSON_56]],\",\n      \"[[PERSON_51]],\"\n    ],\n    [\n      \"[[PERSON_67]],\",\n      \"[[PERSON_57]],\"\n    ],\n    [\n      \"[[PERSON_34]],\",\n      \"[[PERSON_30]],\"\n    ],\n    [\n      \"[[PERSON_16]].\",\n      \"[[PERSON_14]].\"\n    ]\n  ],\n  \"deletions\": [\n    \"[[PERSON_56]] \u2013 \u201eod [[PERSON_56]]\u201c\",\n    \"[[PERSON_57]] \u2013 \u201es [[PERSON_57]]\u201c\",\n    \"[[PERSON_58]] \u2013 \u201eu [[PERSON_58]]\u201c\",\n    \"[[PERSON_59]] \u2013 \u201eo [[PERSON_60]]\u201c\"\n  ],\n  \"insertions\": [\n    \"[[PERSON_21]] \u2013 \u201eu [[PERSON_21]]\u201c\",\n    \"[[PERSON_22]] \u2013 \u201es [[PERSON_22]]\u201c\",\n    \"[[PERSON_23]] \u2013 \u201eo [[PERSON_23]]\u201c\",\n    \"[[PERSON_24]] \u2013 \u201ek [[PERSON_24]]\u201c\"\n  ],\n  \"anchor_insert_before\": \"[[PERSON_25]] \u2013 \u201eu [[PERSON_25]]\u201c\"\n};\n\nconst body = context.document.body;\n\n// 1) Delete the 4 obsolete list paragraphs (by their exact original text).\nfor (const text of DATA.deletions) {\n  const found = body.search(text, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length === 0) {\n    throw new Error(\"Deletion anchor not found: \" + text);\n  }\n  for (const item of found.items) {\n    const para = item.paragraphs.getFirst();\n    para.delete();\n  }\n  await context.sync();\n}\n\n// 2) Insert the 4 new list paragraphs, right before the paragraph that\n//    still reads \"[[PERSON_25]] - \"u [[PERSON_25]]\"\" (anchor unchanged by the diff).\n{\n  const anchorResults = body.search(DATA.anchor_insert_before, { matchCase: true, matchWholeWord: false });\n  anchorResults.load(\"items\");\n  await context.sync();\n  if (anchorResults.items.length === 0) {\n    throw new Error(\"Insertion anchor not found: \" + DATA.anchor_insert_before);\n  }\n  const anchorPara = anchorResults.items[0].paragraphs.getFirst();\n  // Insert in order so the final order matches the diff.\n  let refPara = anchorPara;\n  for (const text of DATA.insertions) {\n    const newPara = refPara.insertParagraph(text, Word.InsertLocation.before);\n    refPara.load(\"text\");\n    await context.sync();\n  }\n}\n\n// 3) Apply the direct text replacements (old run text -> new run text).\nfor (const [oldText, newText] of DATA.replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length === 0) {\n    throw new Error(\"Replacement anchor not found: \" + oldText);\n  }\n  for (const item of found.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Auto-generated edit script for Word COM interop (PowerShell-style).\n# Applies the renumbering of [[PERSON_N]] placeholders across the\n# \"2. SEZNAM SUBJEKTU\" list and all later references in the document,\n# plus inserts 4 new list entries and removes 4 obsolete ones.\n\n$d = $word.ActiveDocument\n\nfunction Find-Paragraph([string]$needle) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $needle\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Anchor not found: $needle\"\n    }\n    $range = $find.Parent\n    return $range.Paragraphs.Item(1)\n}\n\n# 1) Delete the 4 obsolete list paragraphs (by their exact original text).\n$deletions = @(\n    \"[[PERSON_56]] \u2013 \u201eod [[PERSON_56]]\u201c\",\n    \"[[PERSON_57]] \u2013 \u201es [[PERSON_57]]\u201c\",\n    \"[[PERSON_58]] \u2013 \u201eu [[PERSON_58]]\u201c\",\n    \"[[PERSON_59]] \u2013 \u201eo [[PERSON_60]]\u201c\"\n)\nforeach ($needle in $deletions) {\n    $para = Find-Paragraph $needle\n    $para.Range.Delete()\n}\n\n# 2) Insert the 4 new list paragraphs, right before the paragraph that\n#    still reads the anchor text below (unchanged by the diff).\n$anchorText = \"[[PERSON_25]] \u2013 \u201eu [[PERSON_25]]\u201c\"\n$insertions = @(\n    \"[[PERSON_21]] \u2013 \u201eu [[PERSON_21]]\u201c\",\n    \"[[PERSON_22]] \u2013 \u201es [[PERSON_22]]\u201c\",\n    \"[[PERSON_23]] \u2013 \u201eo [[PERSON_23]]\u201c\",\n    \"[[PERSON_24]] \u2013 \u201ek [[PERSON_24]]\u201c\"\n)\n$anchorPara = Find-Paragraph $anchorText\n$insertPoint = $d.Range($anchorPara.Range.Start, $anchorPara.Range.Start)\nforeach ($text in $insertions) {\n    $insertPoint.InsertParagraphBefore()\n    $newRange = $d.Range($insertPoint.Start, $insertPoint.Start)\n    $newRange.Text = $text\n    $insertPoint = $d.Range($anchorPara.Range.Start, $anchorPara.Range.Start)\n}\n\n# 3) Apply the direct text replacements (old run text -> new run text).\n$replacements = @(\n    @{ Old = \"[[PERSON_5]] \u2013 \u201eo [[PERSON_6]]\u201c\"; New = \"[[PERSON_5]] \u2013 \u201eo [[PERSON_5]]\u201c\" },\n    @{ Old = \"[[PERSON_7]] \u2013 \u201ek [[PERSON_8]]\u201c\"; New = \"[[PERSON_6]] \u2013 \u201ek [[PERSON_7]]\u201c\" },\n    @{ Old = \"[[PERSON_9]] \u2013 \u201epro [[PERSON_10]]\u201c\"; New = \"[[PERSON_8]] \u2013 \u201epro [[PERSON_9]]\u201c\" },\n    @{ Old = \"[[PERSON_11]] \u2013 \u201es [[PERSON_12]]\u201c\"; New = \"[[PERSON_10]] \u2013 \u201es [[PERSON_10]]\u201c\" },\n    @{ Old = \"[[PERSON_13]] \u2013 \u201eu [[PERSON_13]]\u201c\"; New = \"[[PERSON_11]] \u2013 \u201eu [[PERSON_11]]\u201c\" },\n    @{ Old = \"[[PERSON_14]] \u2013 \u201eod [[PERSON_14]]\u201c\"; New = \"[[PERSON_12]] \u2013 \u201eod [[PERSON_12]]\u201c\" },\n    @{ Old = \"[[PERSON_15]] \u2013 \u201epro [[PERSON_15]]\u201c\"; New = \"[[PERSON_13]] \u2013 \u201epro [[PERSON_13]]\u201c\" },\n    @{ Old = \"[[PERSON_16]] \u2013 \u201ek [[PERSON_17]]\u201c\"; New = \"[[PERSON_14]] \u2013 \u201ek [[PERSON_15]]\u201c\" },\n    @{ Old = \"[[PERSON_18]] \u2013 \u201es [[PERSON_18]]\u201c\"; New = \"[[PERSON_16]] \u2013 \u201es [[PERSON_16]]\u201c\" },\n    @{ Old = \"[[PERSON_19]] \u2013 \u201eo [[PERSON_20]]\u201c\"; New = \"[[PERSON_17]] \u2013 \u201eo [[PERSON_18]]\u201c\" },\n    @{ Old = \"[[PERSON_21]] \u2013 \u201ek [[PERSON_22]]\u201c\"; New = \"[[PERSON_19]] \u2013 \u201ek [[PERSON_19]]\u201c\" },\n    @{ Old = \"[[PERSON_23]] \u2013 \u201es [[PERSON_24]]\u201c\"; New = \"[[PERSON_20]] \u2013 \u201es [[PERSON_20]]\u201c\" },\n    @{ Old = \"[[PERSON_26]] \u2013 \u201es [[PERSON_26]]\u201c\"; New = \"[[PERSON_26]] \u2013 \u201epro [[PERSON_27]]\u201c\" },\n    @{ Old = \"[[PERSON_27]] \u2013 \u201eo [[PERSON_27]]\u201c\"; New = \"[[PERSON_28]] \u2013 \u201es [[PERSON_28]]\u201c\" },\n    @{ Old = \"[[PERSON_28]] \u2013 \u201ek [[PERSON_28]]\u201c\"; New = \"[[PERSON_29]] \u2013 \u201ek [[PERSON_29]]\u201c\" },\n    @{ Old = \"[[PERSON_29]] \u2013 \u201eu [[PERSON_29]]\u201c\"; New = \"[[PERSON_30]] \u2013 \u201es [[PERSON_31]]\u201c\" },\n    @{ Old = \"[[PERSON_30]] \u2013 \u201epro [[PERSON_31]]\u201c\"; New = \"[[PERSON_32]] \u2013 \u201eo [[PERSON_33]]\u201c\" },\n    @{ Old = \"[[PERSON_32]] \u2013 \u201es [[PERSON_32]]\u201c\"; New = \"[[PERSON_34]] \u2013 \u201epro [[PERSON_34]]\u201c\" },\n    @{ Old = \"[[PERSON_33]] \u2013 \u201ek [[PERSON_33]]\u201c\"; New = \"[[PERSON_35]] \u2013 \u201es [[PERSON_36]]\u201c\" },\n    @{ Old = \"[[PERSON_34]] \u2013 \u201es [[PERSON_35]]\u201c\"; New = \"[[PERSON_37]] \u2013 \u201ek [[PERSON_38]]\u201c\" },\n    @{ Old = \"[[PERSON_36]] \u2013 \u201eo [[PERSON_37]]\u201c\"; New = \"[[PERSON_39]] \u2013 \u201es [[PERSON_39]]\u201c\" },\n    @{ Old = \"[[PERSON_38]] \u2013 \u201epro [[PERSON_38]]\u201c\"; New = \"[[PERSON_40]] \u2013 \u201eo [[PERSON_40]]\u201c\" },\n    @{ Old = \"[[PERSON_39]] \u2013 \u201es [[PERSON_40]]\u201c\"; New = \"[[PERSON_41]] \u2013 \u201eu [[PERSON_41]]\u201c\" },\n    @{ Old = \"[[PERSON_41]] \u2013 \u201ek [[PERSON_42]]\u201c\"; New = \"[[PERSON_42]] \u2013 \u201ek [[PERSON_42]]\u201c\" },\n    @{ Old = \"[[PERSON_43]] \u2013 \u201es [[PERSON_43]]\u201c\"; New = \"[[PERSON_43]] \u2013 \u201ese [[PERSON_44]]\u201c\" },\n    @{ Old = \"[[PERSON_44]] \u2013 \u201eo [[PERSON_45]]\u201c\"; New = \"[[PERSON_45]] \u2013 \u201eu [[PERSON_45]]\u201c\" },\n    @{ Old = \"[[PERSON_46]] \u2013 \u201eu [[PERSON_46]]\u201c\"; New = \"[[PERSON_46]] \u2013 \u201eo [[PERSON_47]]\u201c\" },\n    @{ Old = \"[[PERSON_47]] \u2013 \u201ek [[PERSON_47]]\u201c\"; New = \"[[PERSON_48]] \u2013 \u201es [[PERSON_48]]\u201c\" },\n    @{ Old = \"[[PERSON_48]] \u2013 \u201ese [[PERSON_49]]\u201c\"; New = \"[[PERSON_49]] \u2013 \u201ek [[PERSON_50]]\u201c\" },\n    @{ Old = \"[[PERSON_50]] \u2013 \u201eu [[PERSON_50]]\u201c\"; New = \"[[PERSON_51]] \u2013 \u201eod [[PERSON_51]]\u201c\" },\n    @{ Old = \"[[PERSON_51]] \u2013 \u201eo [[PERSON_52]]\u201c\"; New = \"[[PERSON_52]] \u2013 \u201es [[PERSON_52]]\u201c\" },\n    @{ Old = \"[[PERSON_53]] \u2013 \u201es [[PERSON_53]]\u201c\"; New = \"[[PERSON_53]] \u2013 \u201eu [[PERSON_53]]\u201c\" },\n    @{ Old = \"[[PERSON_54]] \u2013 \u201ek [[PERSON_55]]\u201c\"; New = \"[[PERSON_54]] \u2013 \u201eo [[PERSON_55]]\u201c\" },\n    @{ Old = \"[[PERSON_61]] \u2013 \u201ek [[PERSON_61]]\u201c\"; New = \"[[PERSON_56]] \u2013 \u201ek [[PERSON_56]]\u201c\" },\n    @{ Old = \"V t\u011bchto \u0159\u00edzen\u00edch bylo jedn\u00e1no nap\u0159. s [[PERSON_3]], [[PERSON_9]], [[PERSON_36]] \u010di [[PERSON_62]].\"; New = \"V t\u011bchto \u0159\u00edzen\u00edch bylo jedn\u00e1no nap\u0159. s [[PERSON_3]], [[PERSON_8]], [[PERSON_32]] \u010di [[PERSON_57]].\" },\n    @{ Old = \"sv\u011bdek [[PERSON_43]] (ve v\u00fdpov\u011bdi ozna\u010den jako \u201esv\u011bdek \u010cern\u00e9ho\u201c),\"; New = \"sv\u011bdek [[PERSON_39]] (ve v\u00fdpov\u011bdi ozna\u010den jako \u201esv\u011bdek \u010cern\u00e9ho\u201c),\" },\n    @{ Old = \"po\u0161kozen\u00e1 [[PERSON_21]] (\u201evyj\u00e1d\u0159en\u00ed [[PERSON_21]]\u201c),\"; New = \"po\u0161kozen\u00e1 [[PERSON_19]] (\u201evyj\u00e1d\u0159en\u00ed [[PERSON_19]]\u201c),\" },\n    @{ Old = \"ob\u017ealovan\u00fd [[PERSON_7]] (\u201eobhajoba [[PERSON_7]]\u201c),\"; New = \"ob\u017ealovan\u00fd [[PERSON_6]] (\u201eobhajoba [[PERSON_6]]\u201c),\" },\n    @{ Old = \"pr\u00e1vn\u00ed z\u00e1stupkyn\u011b JUDr. [[PERSON_47]], advok\u00e1tka,\"; New = \"pr\u00e1vn\u00ed z\u00e1stupkyn\u011b JUDr. [[PERSON_42]], advok\u00e1tka,\" },\n    @{ Old = \"tlumo\u010dn\u00edk [[PERSON_46]], zapsan\u00fd v seznamu tlumo\u010dn\u00edk\u016f.\"; New = \"tlumo\u010dn\u00edk [[PERSON_41]], zapsan\u00fd v seznamu tlumo\u010dn\u00edk\u016f.\" },\n    @{ Old = \"Alergologick\u00e9 vy\u0161et\u0159en\u00ed \u010d. ALG/2025/22751 proveden\u00e9 MUDr. [[PERSON_18]],\"; New = \"Alergologick\u00e9 vy\u0161et\u0159en\u00ed \u010d. ALG/2025/22751 proveden\u00e9 MUDr. [[PERSON_16]],\" },\n    @{ Old = \"Neurologick\u00e9 testy \u010d. NEU/2025/44119 proveden\u00e9 MUDr. [[PERSON_51]],\"; New = \"Neurologick\u00e9 testy \u010d. NEU/2025/44119 proveden\u00e9 MUDr. [[PERSON_46]],\" },\n    @{ Old = \"O\u010dn\u00ed vy\u0161et\u0159en\u00ed \u010d. OFT/2023/11281 proveden\u00e9 MUDr. [[PERSON_44]].\"; New = \"O\u010dn\u00ed vy\u0161et\u0159en\u00ed \u010d. OFT/2023/11281 proveden\u00e9 MUDr. [[PERSON_40]].\" },\n    @{ Old = \"Zvl\u00e1\u0161tn\u00ed pozornost byla v\u011bnov\u00e1na v\u00fdsledk\u016fm [[PERSON_25]], [[PERSON_29]] a [[PERSON_59]].\"; New = \"Zvl\u00e1\u0161tn\u00ed pozornost byla v\u011bnov\u00e1na v\u00fdsledk\u016fm [[PERSON_21]], [[PERSON_25]] a [[PERSON_54]].\" },\n    @{ Old = \"mobil [[PERSON_63]] S22, [[IMEI_1]],\"; New = \"mobil [[PERSON_58]] S22, [[IMEI_1]],\" },\n    @{ Old = \"[[PERSON_64]] poskytly technick\u00e9 p\u0159\u00edstupy pro \u0159e\u0161en\u00ed kauz:\"; New = \"[[PERSON_59]] poskytly technick\u00e9 p\u0159\u00edstupy pro \u0159e\u0161en\u00ed kauz:\" },\n    @{ Old = \"pr\u00e1vn\u00ed cloud \u00fa\u010det ID: LEX-ACC-88221 (spravovala [[PERSON_61]]),\"; New = \"pr\u00e1vn\u00ed cloud \u00fa\u010det ID: LEX-ACC-88221 (spravovala [[PERSON_56]]),\" },\n    @{ Old = \"[[PERSON_57]] (\u201ev\u00fdslech [[PERSON_57]]\u201c),\"; New = \"[[PERSON_52]] (\u201ev\u00fdslech [[PERSON_52]]\u201c),\" },\n    @{ Old = \"[[PERSON_53]] (\u201ev\u00fdpov\u011b\u010f [[PERSON_53]]\u201c),\"; New = \"[[PERSON_48]] (\u201ev\u00fdpov\u011b\u010f [[PERSON_48]]\u201c),\" },\n    @{ Old = \"[[PERSON_50]] (\u201ez\u00e1znam o v\u00fdslechu [[PERSON_50]]\u201c),\"; New = \"[[PERSON_45]] (\u201ez\u00e1znam o v\u00fdslechu [[PERSON_45]]\u201c),\" },\n    @{ Old = \"[[PERSON_28]] (\u201ev\u00fdslech [[PERSON_65]]\u201c).\"; New = \"[[PERSON_24]] (\u201ev\u00fdslech [[PERSON_60]]\u201c).\" },\n    @{ Old = \"PhDr. [[PERSON_44]] \u2013 psychologick\u00fd posudek,\"; New = \"PhDr. [[PERSON_40]] \u2013 psychologick\u00fd posudek,\" },\n    @{ Old = \"MUDr. [[PERSON_36]] \u2013 posudek z traumatologie,\"; New = \"MUDr. [[PERSON_32]] \u2013 posudek z traumatologie,\" },\n    @{ Old = \"Ing. [[PERSON_14]] \u2013 expertiza IT infrastruktury.\"; New = \"Ing. [[PERSON_12]] \u2013 expertiza IT infrastruktury.\" },\n    @{ Old = \"Tyto \u00fa\u010dty byly dolo\u017eeny nap\u0159. od [[PERSON_30]], [[PERSON_54]] nebo [[PERSON_66]].\"; New = \"Tyto \u00fa\u010dty byly dolo\u017eeny nap\u0159. od [[PERSON_26]], [[PERSON_49]] nebo [[PERSON_61]].\" },\n    @{ Old = \"[[PERSON_56]],\"; New = \"[[PERSON_51]],\" },\n    @{ Old = \"[[PERSON_67]],\"; New = \"[[PERSON_57]],\" },\n    @{ Old = \"[[PERSON_34]],\"; New = \"[[PERSON_30]],\" },\n    @{ Old = \"[[PERSON_16]].\"; New = \"[[PERSON_14]].\" }\n)\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $ok) {\n        throw \"Replacement anchor not found: $($pair.Old)\"\n    }\n}\n\nWrite-Output \"done\"\n"}
